# Insert a new "Apio" price record for Macroferia Regional de Talca.
# The new row is inserted at row 84, pushing every existing record from
# row 84 onward down by one row (old row 190 becomes row 191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 44740
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112017
$ws.Range("G84").Value = "Apio"
$ws.Range("H84").Value = "Americana (o)"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 800
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 7000
$ws.Range("M84").Value = 7000
$ws.Range("N84").Value = "`$/docena de matas"
$ws.Range("O84").Value = "Provincia del Elquí"
$ws.Range("P84").Value = 1167
$ws.Range("Q84").Value = 6
$ws.Range("R84").Value = "Hortaliza"
